$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Sub Function" status notes for PressureCalibration / syringe Actuation rows
$ws.Range("C30").Value = "More testing"
$ws.Range("C31").Value = "More testing"
$ws.Range("C32").Value = "More testing"

# Reflect the view/selection state at save time (scrolled down, C30 selected)
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("C30").Select() | Out-Null
